$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 43998.914
$ws.Range("J33").Value = 1021.55554
$ws.Range("L33").Value = 1021.55554
$ws.Range("N33").Value = -1479.55554
$ws.Range("H38").Value = 1454
$ws.Range("I38").Value = 124.25
$ws.Range("J38").Value = 5000
$ws.Range("K38").Value = 372.75
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = -0.75
$ws.Range("N38").Value = -15744
$ws.Range("H39").Value = 1065.8334
$ws.Range("I39").Value = 890
$ws.Range("J39").Value = 3000
$ws.Range("K39").Value = 2670
$ws.Range("L39").Value = 9000
$ws.Range("M39").Value = -2374
$ws.Range("N39").Value = -9592
$ws.Range("H40").Value = 8697.777
$ws.Range("J40").Value = 4111.4287
$ws.Range("L40").Value = 4111.4287
$ws.Range("N40").Value = -4461.4287
$ws.Range("H41").Value = 3256.375
$ws.Range("I41").Value = 3310.8
$ws.Range("K41").Value = 3310.8
$ws.Range("M41").Value = -2870.8
$ws.Range("H53").Value = 205.4
$ws.Range("I53").Value = 118.25
$ws.Range("J53").Value = 263.5
$ws.Range("K53").Value = 118.25
$ws.Range("L53").Value = 263.5
$ws.Range("M53").Value = 518.75
$ws.Range("N53").Value = -1537.5
$ws.Range("H55").Value = 1086.2667
$ws.Range("I55").Value = 464.2857
$ws.Range("K55").Value = 464.2857
$ws.Range("M55").Value = -250.2857
$ws.Range("H64").Value = 5815.4614
$ws.Range("I64").Value = 5781.909
$ws.Range("K64").Value = 5781.909
$ws.Range("M64").Value = -5533.909
$ws.Range("H67").Value = 5815.4614
$ws.Range("I67").Value = 5781.909
$ws.Range("K67").Value = 5781.909
$ws.Range("M67").Value = -4923.909
$ws.Range("H92").Value = 747.8889
$ws.Range("J92").Value = 1300
$ws.Range("L92").Value = 1300
$ws.Range("N92").Value = -3796
$ws.Range("H94").Value = 1219
$ws.Range("I94").Value = 1219
$ws.Range("K94").Value = 1219
$ws.Range("M94").Value = -768
$ws.Range("H96").Value = 1317.5264
$ws.Range("I96").Value = 1257.2727
$ws.Range("J96").Value = 1400.375
$ws.Range("K96").Value = 3771.8181
$ws.Range("L96").Value = 4201.125
$ws.Range("M96").Value = -2398.8181
$ws.Range("N96").Value = -6947.125
$ws.Range("H99").Value = 1060.6666
$ws.Range("J99").Value = 2250
$ws.Range("L99").Value = 6750
$ws.Range("N99").Value = -9746
$ws.Range("H100").Value = 8133.7
$ws.Range("I100").Value = 8762.571
$ws.Range("J100").Value = 6666.3335
$ws.Range("K100").Value = 8762.571
$ws.Range("L100").Value = 6666.3335
$ws.Range("M100").Value = -8221.571
$ws.Range("N100").Value = -7748.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4700.8486
$ws.Range("I32").Value = 4700.8486
$ws.Range("K32").Value = 4700.8486
$ws.Range("M32").Value = -4413.8486
$ws.Range("H61").Value = 1996.75
$ws.Range("I61").Value = 1996.75
$ws.Range("K61").Value = 1996.75
$ws.Range("M61").Value = -1784.75
$ws.Range("H62").Value = 89405.836
$ws.Range("J62").Value = 89405.836
$ws.Range("L62").Value = 89405.836
$ws.Range("N62").Value = -90653.836
$ws.Range("H65").Value = 89405.836
$ws.Range("J65").Value = 89405.836
$ws.Range("L65").Value = 268217.508
$ws.Range("N65").Value = -274457.508
$ws.Range("H75").Value = 43299.668
$ws.Range("J75").Value = 43299.668
$ws.Range("L75").Value = 43299.668
$ws.Range("N75").Value = -45047.668
$ws.Range("H78").Value = 43299.668
$ws.Range("J78").Value = 43299.668
$ws.Range("L78").Value = 129899.004
$ws.Range("N78").Value = -138635.004
$ws.Range("H88").Value = 1177.1578
$ws.Range("I88").Value = 1026
$ws.Range("J88").Value = 1313.2
$ws.Range("K88").Value = 1026
$ws.Range("L88").Value = 1313.2
$ws.Range("M88").Value = -620
$ws.Range("N88").Value = -2125.2
$ws.Range("H91").Value = 1177.1578
$ws.Range("I91").Value = 1026
$ws.Range("J91").Value = 1313.2
$ws.Range("K91").Value = 1026
$ws.Range("L91").Value = 1313.2
$ws.Range("M91").Value = 378
$ws.Range("N91").Value = -4121.2
$ws.Range("H136").Value = 1996.75
$ws.Range("I136").Value = 1996.75
$ws.Range("K136").Value = 5990.25
$ws.Range("M136").Value = -3440.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 391.91666
$ws.Range("J80").Value = 570.8570999999999
$ws.Range("L80").Value = 570.8570999999999
$ws.Range("N80").Value = -2566.8571
$ws.Range("H83").Value = 391.91666
$ws.Range("J83").Value = 570.8570999999999
$ws.Range("L83").Value = 2854.2855
$ws.Range("N83").Value = -12838.2855
$ws.Range("H94").Value = 3147.75
$ws.Range("J94").Value = 1546.6666
$ws.Range("L94").Value = 1546.6666
$ws.Range("N94").Value = -2448.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 22333.5
$ws.Range("J50").Value = 22333.5
$ws.Range("L50").Value = 22333.5
$ws.Range("N50").Value = -23583.5
$ws.Range("H58").Value = 2559.9092
$ws.Range("I58").Value = 2122.625
$ws.Range("K58").Value = 2122.625
$ws.Range("M58").Value = -1919.625
$ws.Range("H59").Value = 23119.334
$ws.Range("J59").Value = 25127
$ws.Range("L59").Value = 25127
$ws.Range("N59").Value = -27417
$ws.Range("H60").Value = 11110.777
$ws.Range("J60").Value = 19001.666
$ws.Range("L60").Value = 19001.666
$ws.Range("N60").Value = -20023.666
$ws.Range("H136").Value = 2559.9092
$ws.Range("I136").Value = 2122.625
$ws.Range("K136").Value = 6367.875
$ws.Range("M136").Value = -3817.875
$ws.Range("H141").Value = 180624.88
$ws.Range("J141").Value = 180624.88
$ws.Range("L141").Value = 180624.88
$ws.Range("N141").Value = -190984.88

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 44.727272
$ws.Range("I2").Value = 34.57143
$ws.Range("K2").Value = 207.42858
$ws.Range("M2").Value = -94.42858000000001
$ws.Range("H12").Value = 261.8
$ws.Range("I12").Value = 108.26667
$ws.Range("K12").Value = 324.80001
$ws.Range("M12").Value = -151.80001
$ws.Range("H33").Value = 39
$ws.Range("J33").Value = 30
$ws.Range("L33").Value = 180
$ws.Range("N33").Value = -746
$ws.Range("H68").Value = 2023.625
$ws.Range("I68").Value = 2066.3333
$ws.Range("J68").Value = 1998
$ws.Range("K68").Value = 6198.999899999999
$ws.Range("L68").Value = 5994
$ws.Range("M68").Value = -5387.999899999999
$ws.Range("N68").Value = -7616
$ws.Range("H71").Value = 2023.625
$ws.Range("I71").Value = 2066.3333
$ws.Range("J71").Value = 1998
$ws.Range("K71").Value = 18596.9997
$ws.Range("L71").Value = 17982
$ws.Range("M71").Value = -14540.9997
$ws.Range("N71").Value = -26094
$ws.Range("H98").Value = 1570.8572
$ws.Range("I98").Value = 1599.6
$ws.Range("K98").Value = 4798.799999999999
$ws.Range("M98").Value = -3300.799999999999
$ws.Range("H104").Value = 2848.5833
$ws.Range("J104").Value = 3132.5557
$ws.Range("L104").Value = 9397.667099999999
$ws.Range("N104").Value = -14639.6671
$ws.Range("H130").Value = 4599.7144
$ws.Range("J130").Value = 5033
$ws.Range("L130").Value = 15099
$ws.Range("N130").Value = -25139

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 28000
$ws.Range("J70").Value = 28000
$ws.Range("L70").Value = 28000
$ws.Range("N70").Value = -28540
$ws.Range("H73").Value = 28000
$ws.Range("J73").Value = 28000
$ws.Range("L73").Value = 28000
$ws.Range("N73").Value = -29872
$ws.Range("H100").Value = 601635.8
$ws.Range("I100").Value = 751419.75
$ws.Range("K100").Value = 751419.75
$ws.Range("M100").Value = -750878.75
$ws.Range("H122").Value = 2717
$ws.Range("I122").Value = 2208.8333
$ws.Range("K122").Value = 6626.499899999999
$ws.Range("M122").Value = -4176.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 20531.666
$ws.Range("J11").Value = 28297.5
$ws.Range("L11").Value = 28297.5
$ws.Range("N11").Value = -28581.5
$ws.Range("H12").Value = 50006
$ws.Range("I12").Value = 50006
$ws.Range("K12").Value = 50006
$ws.Range("M12").Value = -49864
$ws.Range("H30").Value = 24980
$ws.Range("J30").Value = 24980
$ws.Range("L30").Value = 24980
$ws.Range("N30").Value = -25194
$ws.Range("H96").Value = 2718.8
$ws.Range("I96").Value = 3901.3333
$ws.Range("J96").Value = 945
$ws.Range("K96").Value = 3901.3333
$ws.Range("L96").Value = 945
$ws.Range("M96").Value = -2528.3333
$ws.Range("N96").Value = -3691
$ws.Range("H119").Value = 19505
$ws.Range("J119").Value = 19505
$ws.Range("L119").Value = 19505
$ws.Range("N119").Value = -27680
